# "upd. order of the extremums"
# Reorders the data rows (rows 5-18, columns A:T) of the single worksheet:
# each row's full content (values + number types + styles) is relocated to a
# new row position per the mapping below. Row 13 stays in place.
#
# new row  <=  old row (source of the moved content)
#    5     <=   12
#    6     <=   14
#    7     <=   17
#    8     <=    5
#    9     <=    6
#   10     <=    7
#   11     <=   10
#   12     <=   11
#   13     <=   13   (unchanged)
#   14     <=    8
#   15     <=    9
#   16     <=   15
#   17     <=   16
#   18     <=   18   (unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage the original rows 5-18 off to the side first (far enough away that it
# can't collide with any destination), so every destination write reads from
# an untouched copy of the original data, regardless of write order.
$stagingFirstRow = 100
for ($r = 5; $r -le 18; $r++) {
    $srcRow = $r
    $stageRow = $stagingFirstRow + ($r - 5)
    $ws.Range("A" + $srcRow + ":T" + $srcRow).Copy($ws.Range("A" + $stageRow + ":T" + $stageRow))
}

# old row -> new row mapping (expressed the other way round from the staged copies)
$moves = @{
    12 = 5
    14 = 6
    17 = 7
    5  = 8
    6  = 9
    7  = 10
    10 = 11
    11 = 12
    13 = 13
    8  = 14
    9  = 15
    15 = 16
    16 = 17
    18 = 18
}

foreach ($oldRow in $moves.Keys) {
    $newRow = $moves[$oldRow]
    $stageRow = $stagingFirstRow + ($oldRow - 5)
    $ws.Range("A" + $stageRow + ":T" + $stageRow).Copy($ws.Range("A" + $newRow + ":T" + $newRow))
}

# Clear the staging area again so it doesn't leave stray data behind.
$ws.Range("A" + $stagingFirstRow + ":T" + ($stagingFirstRow + 13)).Clear()

# The saved selection in the sheet view also moved (B26 -> B24).
[void]$ws.Range("B24").Select()
